$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.818.99"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "1.649.28"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.76"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0629"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.37"
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").Value = "1.879.01"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").Value = "1.666.65"
$ws.Range("E13").Value = "  +2.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  +1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  +1.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.95"
$ws.Range("E16").Value = "  +4.48%  "

$ws.Range("D17").Value = "26.838.51"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.30"
$ws.Range("E19").Value = "  +4.37%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +2.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.49"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +7.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.85"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.83"
$ws.Range("E29").Value = "  +2.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0521"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +3.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +2.45%  "

$ws.Range("D34").Value = "1.277.99"
$ws.Range("E34").Value = "  +9.32%  "

$ws.Range("E35").Value = "  +1.59%  "

$ws.Range("E36").Value = "  +1.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  +4.34%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.516"
$ws.Range("E38").Value = "  +2.38%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.810"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("D44").Value = "1.788.33"
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.81"
$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("E46").Value = "  +3.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.03"
$ws.Range("E47").Value = "  +2.68%  "

$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.71"
$ws.Range("E50").Value = "  +2.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0967"
$ws.Range("E51").Value = "  +3.08%  "
